$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows before row 32, pushing old rows 32-36 down to 37-41
$ws.Range("A32:A36").EntireRow.Insert()

# Fill in the new rows 32-36 with the new weekly data (week of 2021-09-24 / serial 44463)
$newRows = @(
    @{ Row=32; D=44463; H="Sin especificar"; I="Banquete"; J=34;   K=23000; L=23000; M=23000; N="$/bandeja 10 kilos"; O="Provincia de Linares"; P=2300; Q=10 },
    @{ Row=33; D=44463; H="Sin especificar"; I="Primera";  J=52;   K=21000; L=21000; M=21000; N="$/bandeja 10 kilos"; O="Provincia de Linares"; P=2100; Q=10 },
    @{ Row=34; D=44463; H="Sin especificar"; I="Primera";  J=340;  K=2000;  L=2000;  M=2000;  N="$/kilo";             O="Provincia de Linares"; P=2000; Q=1 },
    @{ Row=35; D=44463; H="Sin especificar"; I="Segunda";  J=43;   K=18000; L=18000; M=18000; N="$/bandeja 10 kilos"; O="Provincia de Linares"; P=1800; Q=10 },
    @{ Row=36; D=44463; H="Sin especificar"; I="Segunda";  J=160;  K=1700;  L=1700;  M=1700;  N="$/kilo";             O="Provincia de Linares"; P=1700; Q=1 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = 9
    $ws.Cells.Item($row, 2).Value2 = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($row, 3).Value2 = "Metropolitana"
    $ws.Cells.Item($row, 4).Value2 = $r.D
    $ws.Cells.Item($row, 5).Value2 = 13
    $ws.Cells.Item($row, 6).Value2 = 300000000
    $ws.Cells.Item($row, 7).Value2 = "Espárragos"
    $ws.Cells.Item($row, 8).Value2 = $r.H
    $ws.Cells.Item($row, 9).Value2 = $r.I
    $ws.Cells.Item($row, 10).Value2 = $r.J
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $r.N
    $ws.Cells.Item($row, 15).Value2 = $r.O
    $ws.Cells.Item($row, 16).Value2 = $r.P
    $ws.Cells.Item($row, 17).Value2 = $r.Q
    $ws.Cells.Item($row, 18).Value2 = "Hortaliza"
}
